$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table with the latest scraped snapshot.
# A few rows also had their coin reshuffled (rank swap), so Coin/Link cells
# are rewritten too in those spots.
#
# Several Price (column D) values are plain decimal-looking text
# (e.g. "12.80", "1.00") that originally round-trip as literal strings
# (inline/shared-string cells), not numbers. Excel's smart Value setter
# would otherwise coerce them to numbers and normalize away trailing
# zeros (e.g. "12.80" -> 12.8), so NumberFormat is forced to Text ("@")
# immediately before writing those specific cells to keep them literal.

$ws.Range("D2").Value = "69.847.54"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.525.45"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.52"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.88"
$ws.Range("E6").Value = "  +6.14%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -6.85%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.82"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "4.084.78"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "598.79"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.80"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.010.51"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.10"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "3.525.74"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.23"
$ws.Range("E22").Value = "  +6.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.30"
$ws.Range("E23").Value = "  +5.03%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.68"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.30"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.18"
$ws.Range("E26").Value = "  +5.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.62"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +11.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.44"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.12"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.0₃0852"
$ws.Range("E35").Value = "  +10.10%  "
$ws.Range("D36").Value = "3.735.91"
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.63"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "487.69"
$ws.Range("E42").Value = "  -5.99%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  +11.38%  "
